$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold price text that looks numeric (e.g. "214.30", "4.06").
# Excel would otherwise coerce these into floating point numbers and lose the
# exact textual representation (trailing zeros / float rounding). Forcing the
# cell to Text format before the write keeps the value an exact string, and
# resetting the Style back to Normal afterwards avoids leaving a stray
# number-format style on the cell (matching the original unformatted cells).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.816.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.85%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.637.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.13%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("E8").Value = "  +1.36%  "
$ws.Range("E9").Value = "  +0.66%  "
$ws.Range("E10").Value = "  +2.48%  "
$ws.Range("E11").Value = "  +0.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.867.04"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.636.01"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.515"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.90%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "241.71"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.797.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0727"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.56%  "
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.23%  "
$ws.Range("E23").Value = "  +2.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("E27").Value = "  +1.79%  "
$ws.Range("E28").Value = "  +0.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.83"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0497"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.47%  "
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.511.51"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.47%  "
$ws.Range("E33").Value = "  +1.93%  "
$ws.Range("E34").Value = "  +2.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.56"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.83%  "
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.576"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.95%  "
$ws.Range("E38").Value = "  +1.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.858"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.94"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.95%  "
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("E42").Value = "  +1.27%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "63.80"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.69%  "
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.777.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.766"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.911"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.53"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0976"
$ws.Range("D49").Style = "Normal"
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.28%  "
